# Update the "Prix Spot" sheet: insert a new date column (20-nov) before the
# 01-oct. column block, shifting DU:EY -> DV:EZ.
$wb = $excel.ActiveWorkbook

$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Columns("DU").Insert()
$wsPrix.Range("DU1").Value = "20-nov"
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 125).Value = "-"
}

# Append the latest daily price row to "Gaz".
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A154").NumberFormat = "@"
$wsGaz.Range("A154").Value = "2025-11-18"
$wsGaz.Range("A154").ClearFormats()
$wsGaz.Range("B154").Value = 30.615

# Append the latest daily price row to "CO2".
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A154").NumberFormat = "@"
$wsCO2.Range("A154").Value = "2025-11-18"
$wsCO2.Range("A154").ClearFormats()
$wsCO2.Range("B154").Value = 80.93000000000001
